# Word COM-interop script that applies the edit described by the diff:
#  1. Appends a red-colored "(This is a change – Version for branch
#     alternate)" annotation, split across three runs, to the first
#     paragraph ("This is a Microsoft word document."), after two
#     trailing spaces added to the original sentence.
#  2. Adds a new, empty paragraph (shaded background fill F9F9F9) right
#     after the final "Free at last..." paragraph, just before the
#     section break.

$d = $word.ActiveDocument

# --- 1. Extend the first paragraph -----------------------------------

$firstParaLen = $d.Paragraphs(1).Range.Text.Length - 1   # drop trailing pilcrow

# Two plain trailing spaces (no special formatting) appended directly
# after the existing sentence.
$r1 = $d.Range($firstParaLen, $firstParaLen)
$r1.InsertAfter("  ")

$dash = [char]0x2013   # EN DASH ("–")

# Run 2: "(This is a change – Ve"  (dark red, C00000)
$pos2 = $firstParaLen + 2
$chunk2 = "(This is a change " + $dash + " Ve"
$r2 = $d.Range($pos2, $pos2)
$r2.InsertAfter($chunk2)
$r2.Font.Color = 192   # RGB(192,0,0) -> w:color w:val="C00000"

# Run 3: "rsion for branch alternate" (dark red, C00000)
$pos3 = $pos2 + $chunk2.Length
$chunk3 = "rsion for branch alternate"
$r3 = $d.Range($pos3, $pos3)
$r3.InsertAfter($chunk3)
$r3.Font.Color = 192

# Run 4: ")" (dark red, C00000)
$pos4 = $pos3 + $chunk3.Length
$chunk4 = ")"
$r4 = $d.Range($pos4, $pos4)
$r4.InsertAfter($chunk4)
$r4.Font.Color = 192

# --- 2. Append a new shaded, empty paragraph at the very end ---------

$tail = $d.Content
$tail.Collapse(0)   # wdCollapseEnd
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>'
$tail.InsertXML($newParaXml)
